$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.832.38"
$ws.Range("E2").Value = "  -8.98%  "

$ws.Range("D3").Value = "2.421.74"
$ws.Range("E3").Value = "  -11.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "462.30"
$ws.Range("E5").Value = "  -9.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.88"
$ws.Range("E6").Value = "  -8.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.484"
$ws.Range("E8").Value = "  -9.65%  "

$ws.Range("D9").Value = "2.412.30"
$ws.Range("E9").Value = "  -12.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0946"
$ws.Range("E10").Value = "  -10.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.31"
$ws.Range("E11").Value = "  -12.81%  "

$ws.Range("E12").Value = "  -10.97%  "

$ws.Range("E13").Value = "  -4.42%  "

$ws.Range("D14").Value = "2.816.10"
$ws.Range("E14").Value = "  -12.36%  "

$ws.Range("D15").Value = "53.738.18"
$ws.Range("E15").Value = "  -9.17%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000132"
$ws.Range("E16").Value = "  -3.70%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.60"
$ws.Range("E17").Value = "  -10.45%  "

$ws.Range("D18").Value = "2.412.75"
$ws.Range("E18").Value = "  -12.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.16"
$ws.Range("E19").Value = "  -13.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "306.78"
$ws.Range("E20").Value = "  -11.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.36"
$ws.Range("E21").Value = "  -15.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.33"
$ws.Range("E24").Value = "  -15.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "55.82"
$ws.Range("E25").Value = "  -12.18%  "

$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("E27").Value = "  -10.67%  "

$ws.Range("E28").Value = "  -12.07%  "

$ws.Range("D29").Value = "2.495.21"
$ws.Range("E29").Value = "  -12.82%  "

$ws.Range("E30").Value = "  -6.60%  "

$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").Value = "0.0₃0716"
$ws.Range("E32").Value = "  -15.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "145.70"
$ws.Range("E33").Value = "  -2.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.63"
$ws.Range("E34").Value = "  -8.81%  "

$ws.Range("E35").Value = "  -12.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.97"
$ws.Range("E36").Value = "  -8.78%  "

$ws.Range("E37").Value = "  -17.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.05"
$ws.Range("E38").Value = "  -7.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.800"
$ws.Range("E39").Value = "  -16.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "32.73"
$ws.Range("E41").Value = "  -9.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.590"
$ws.Range("E42").Value = "  -3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0522"
$ws.Range("E43").Value = "  -6.95%  "

$ws.Range("E44").Value = "  -8.82%  "

$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("E46").Value = "  -12.16%  "

$ws.Range("D47").Value = "1.928.01"
$ws.Range("E47").Value = "  -12.17%  "

$ws.Range("E48").Value = "  -4.48%  "

$ws.Range("E49").Value = "  -2.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.18"
$ws.Range("E50").Value = "  -12.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.42"
$ws.Range("E51").Value = "  -14.47%  "
